# Updates the cryptos worksheet with refreshed prices / 1h volume figures
# (values sourced from the scheduled "Updated cryptos list" GitHub Actions run).
# A handful of rows also had their coin order swapped (Hedera / WEMIXToken).
#
# Price and volume columns are stored as plain text in the workbook (e.g.
# "39.967.23" or "  -4.25%  "), not numbers, so that formatting such as
# thousand-separating dots, leading "0.0"-style precision and padding
# spaces survives untouched. Excel will happily re-interpret a numeric-
# looking string (like "0.530") as a real number and silently drop the
# trailing zero, so for any replacement value that looks like a plain
# number we force the cell to Text format first ("@") before writing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '40.114.42'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -3.95%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '2.338.90'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -5.63%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  -0.13%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '308.33'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -3.98%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '85.13'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -7.80%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.530'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -3.82%  '; ForceText = $false },
    @{ Cell = 'E8'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.486'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -5.16%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '0.0823'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  -4.59%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '30.32'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -8.41%  '; ForceText = $false },
    @{ Cell = 'E12'; Value = '  -0.32%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '2.695.87'; ForceText = $false },
    @{ Cell = 'E13'; Value = '  -5.76%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '6.45'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -6.59%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '14.78'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -4.93%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '2.329.92'; ForceText = $false },
    @{ Cell = 'E16'; Value = '  -5.77%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '0.752'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -5.52%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '40.013.39'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  -4.06%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '0.0₃0907'; ForceText = $false },
    @{ Cell = 'E19'; Value = '  -3.88%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '6.10'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  -5.32%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '67.80'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -4.25%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '10.72'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -4.91%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '235.76'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -1.78%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '2.55'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -7.39%  '; ForceText = $false },
    @{ Cell = 'E25'; Value = '  +0.18%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '1.81'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -6.97%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '23.53'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -6.04%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '2.22'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -1.02%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '9.22'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -5.43%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '35.12'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -4.35%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '151.50'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -3.76%  '; ForceText = $false },
    @{ Cell = 'E32'; Value = '  -0.17%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '5.14'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -5.68%  '; ForceText = $false },
    @{ Cell = 'B34'; Value = 'WEMIXToken'; ForceText = $false },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; ForceText = $false },
    @{ Cell = 'D34'; Value = '2.45'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -4.65%  '; ForceText = $false },
    @{ Cell = 'B35'; Value = 'Hedera'; ForceText = $false },
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false },
    @{ Cell = 'D35'; Value = '0.0726'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  -5.09%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '0.114'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -2.13%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '2.78'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  -3.35%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.100'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -3.45%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '15.85'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -7.89%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '1.71'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -7.39%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '3.82'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -4.81%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '2.31'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -4.56%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '1.947.55'; ForceText = $false },
    @{ Cell = 'E43'; Value = '  -2.49%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '0.0267'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -6.08%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '17.82'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -4.96%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '9.31'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -1.57%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '2.69'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -9.70%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '2.553.35'; ForceText = $false },
    @{ Cell = 'E48'; Value = '  -6.72%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '92.89'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -5.00%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '71.17'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -6.70%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '63.64'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -5.64%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
